# This workbook contains a weekly/daily "Apio" (celery) price log for the
# "Vega Modelo de Temuco" market. The update adds one new daily record,
# which is inserted as a new row 377, pushing all the subsequent rows
# (previously 377-493) down by one (to 378-494).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 377, shifting existing rows 377-493 down to 378-494.
$ws.Rows.Item(377).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A377").Value = 10
$ws.Range("B377").Value = "Vega Modelo de Temuco"
$ws.Range("C377").Value = "La Araucanía"
$ws.Range("D377").Value = 45093
$ws.Range("E377").Value = 9
$ws.Range("F377").Value = 100112017
$ws.Range("G377").Value = "Apio"
$ws.Range("H377").Value = "Americana (o)"
$ws.Range("I377").Value = "Primera"
$ws.Range("J377").Value = 65
$ws.Range("K377").Value = 8000
$ws.Range("L377").Value = 8000
$ws.Range("M377").Value = 8000
$ws.Range("N377").Value = "`$/docena de matas"
$ws.Range("O377").Value = "Provincia del Elquí"
$ws.Range("P377").Value = 1333
$ws.Range("Q377").Value = 6
$ws.Range("R377").Value = "Hortaliza"
